$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New season-record columns: Wins (AD), Losses (AE), Ties (AF)
$lastRow = 58

# Header row (row 1) - labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, border, centered) from the
# neighboring existing header cell (AC1) onto the new header cells.
$headerFormatSrc = $ws.Range("AC1")
$headerFormatDst = $ws.Range("AD1:AF1")
$headerFormatSrc.Copy()
$headerFormatDst.PasteSpecial(-4122)

# Fill in the season record for every data row (2015 NYY: 87-75-0)
$ws.Range("AD2:AD" + $lastRow).Value = 87
$ws.Range("AE2:AE" + $lastRow).Value = 75
$ws.Range("AF2:AF" + $lastRow).Value = 0
